$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (rows 2-28) holds a "Förändrad" (last-changed) date stored as a
# raw Excel serial date number. Bump it from 45430 (2024-05-18) to
# 45431 (2024-05-19) for every data row.
for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45430) {
        $cell.Value2 = 45431
    }
}
